$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header cells for the new Wins/Losses/Ties columns, matching the
# existing header style (s="1") used by the other header cells in row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill in the team record (Wins=62, Losses=100, Ties=0) for every data row.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 62
    $ws.Cells.Item($r, 31).Value = 100
    $ws.Cells.Item($r, 32).Value = 0
}
